$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

if ($t.Rows.Count -ne 20 -or $t.Columns.Count -ne 5) {
    throw ("Unexpected table shape: " + $t.Rows.Count + "x" + $t.Columns.Count)
}

$t.Cell(1,1).Range.Text = "19+41=60"
$t.Cell(1,2).Range.Text = "6+54=60"
$t.Cell(1,3).Range.Text = "48+26=74"
$t.Cell(1,4).Range.Text = "35+33=68"
$t.Cell(1,5).Range.Text = "54+12=66"
$t.Cell(2,1).Range.Text = "73+24=97"
$t.Cell(2,2).Range.Text = "71-57=14"
$t.Cell(2,3).Range.Text = "58+6=64"
$t.Cell(2,4).Range.Text = "65-51=14"
$t.Cell(2,5).Range.Text = "97-28=69"
$t.Cell(3,1).Range.Text = "71+13=84"
$t.Cell(3,2).Range.Text = "5+78=83"
$t.Cell(3,3).Range.Text = "78-29=49"
$t.Cell(3,4).Range.Text = "11-4=7"
$t.Cell(3,5).Range.Text = "55+35=90"
$t.Cell(4,1).Range.Text = "13-5=8"
$t.Cell(4,2).Range.Text = "58+27=85"
$t.Cell(4,3).Range.Text = "24+56=80"
$t.Cell(4,4).Range.Text = "9-7=2"
$t.Cell(4,5).Range.Text = "23+73=96"
$t.Cell(5,1).Range.Text = "2+72=74"
$t.Cell(5,2).Range.Text = "53-16=37"
$t.Cell(5,3).Range.Text = "6+58=64"
$t.Cell(5,4).Range.Text = "14+32=46"
$t.Cell(5,5).Range.Text = "75-25=50"
$t.Cell(6,1).Range.Text = "62-21=41"
$t.Cell(6,2).Range.Text = "55+1=56"
$t.Cell(6,3).Range.Text = "21-3=18"
$t.Cell(6,4).Range.Text = "9+32=41"
$t.Cell(6,5).Range.Text = "94-59=35"
$t.Cell(7,1).Range.Text = "50+46=96"
$t.Cell(7,2).Range.Text = "83-57=26"
$t.Cell(7,3).Range.Text = "37+28=65"
$t.Cell(7,4).Range.Text = "94-89=5"
$t.Cell(7,5).Range.Text = "7+17=24"
$t.Cell(8,1).Range.Text = "77-48=29"
$t.Cell(8,2).Range.Text = "78-75=3"
$t.Cell(8,3).Range.Text = "73+24=97"
$t.Cell(8,4).Range.Text = "2+18=20"
$t.Cell(8,5).Range.Text = "63-60=3"
$t.Cell(9,1).Range.Text = "11+66=77"
$t.Cell(9,2).Range.Text = "35+32=67"
$t.Cell(9,3).Range.Text = "8+65=73"
$t.Cell(9,4).Range.Text = "80-39=41"
$t.Cell(9,5).Range.Text = "49-27=22"
$t.Cell(10,1).Range.Text = "98-67=31"
$t.Cell(10,2).Range.Text = "29+44=73"
$t.Cell(10,3).Range.Text = "3+40=43"
$t.Cell(10,4).Range.Text = "39+6=45"
$t.Cell(10,5).Range.Text = "93-84=9"
$t.Cell(11,1).Range.Text = "55+16=71"
$t.Cell(11,2).Range.Text = "37-27=10"
$t.Cell(11,3).Range.Text = "64+21=85"
$t.Cell(11,4).Range.Text = "53-36=17"
$t.Cell(11,5).Range.Text = "21+45=66"
$t.Cell(12,1).Range.Text = "89-14=75"
$t.Cell(12,2).Range.Text = "44-28=16"
$t.Cell(12,3).Range.Text = "71-49=22"
$t.Cell(12,4).Range.Text = "64-14=50"
$t.Cell(12,5).Range.Text = "51+1=52"
$t.Cell(13,1).Range.Text = "94-61=33"
$t.Cell(13,2).Range.Text = "60+16=76"
$t.Cell(13,3).Range.Text = "26+15=41"
$t.Cell(13,4).Range.Text = "39+14=53"
$t.Cell(13,5).Range.Text = "29-10=19"
$t.Cell(14,1).Range.Text = "80-64=16"
$t.Cell(14,2).Range.Text = "62-37=25"
$t.Cell(14,3).Range.Text = "6+35=41"
$t.Cell(14,4).Range.Text = "85-7=78"
$t.Cell(14,5).Range.Text = "7+46=53"
$t.Cell(15,1).Range.Text = "49-3=46"
$t.Cell(15,2).Range.Text = "64-46=18"
$t.Cell(15,3).Range.Text = "5+35=40"
$t.Cell(15,4).Range.Text = "89-57=32"
$t.Cell(15,5).Range.Text = "34-4=30"
$t.Cell(16,1).Range.Text = "37-7=30"
$t.Cell(16,2).Range.Text = "53+1=54"
$t.Cell(16,3).Range.Text = "14+76=90"
$t.Cell(16,4).Range.Text = "20-13=7"
$t.Cell(16,5).Range.Text = "19+29=48"
$t.Cell(17,1).Range.Text = "59+35=94"
$t.Cell(17,2).Range.Text = "57-19=38"
$t.Cell(17,3).Range.Text = "96-50=46"
$t.Cell(17,4).Range.Text = "30+41=71"
$t.Cell(17,5).Range.Text = "61+3=64"
$t.Cell(18,1).Range.Text = "32-25=7"
$t.Cell(18,2).Range.Text = "60-4=56"
$t.Cell(18,3).Range.Text = "48+15=63"
$t.Cell(18,4).Range.Text = "57+28=85"
$t.Cell(18,5).Range.Text = "87-55=32"
$t.Cell(19,1).Range.Text = "70-31=39"
$t.Cell(19,2).Range.Text = "32-3=29"
$t.Cell(19,3).Range.Text = "11+40=51"
$t.Cell(19,4).Range.Text = "2+27=29"
$t.Cell(19,5).Range.Text = "93-2=91"
$t.Cell(20,1).Range.Text = "88-10=78"
$t.Cell(20,2).Range.Text = "59-37=22"
$t.Cell(20,3).Range.Text = "98-29=69"
$t.Cell(20,4).Range.Text = "38-18=20"
$t.Cell(20,5).Range.Text = "56+36=92"

Write-Output ("Done. Rows=" + $t.Rows.Count + " Cols=" + $t.Columns.Count)
